$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
# A8: "Volume 32   Number  14" -> "Volume 32   Number  15"
$ws.Range("A8").Value = "Volume 32   Number  15"
# C9: "Report Covering the Week  3/31/2025  Through  4/6/2025"
#  -> "Report Covering the Week  4/7/2025  Through  4/13/2025"
$ws.Range("C9").Value = "Report Covering the Week  4/7/2025  Through  4/13/2025"

# --- Crime-complaint table updates (rows 15-33) ---
# Row 15
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 7
$ws.Range("H15").Value = -85.714285714285
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = -46.153846153846
$ws.Range("L15").Value = 40

# Row 16
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -33.333333333333
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -22.727272727272
$ws.Range("I16").Value = 70
$ws.Range("J16").Value = 80
$ws.Range("K16").Value = -12.5
$ws.Range("L16").Value = 4.477611940298
$ws.Range("M16").Value = -12.5
$ws.Range("N16").Value = -67.441860465116

# Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = 3.703703703703
$ws.Range("I17").Value = 112
$ws.Range("J17").Value = 114
$ws.Range("K17").Value = -1.754385964912
$ws.Range("L17").Value = 3.703703703703
$ws.Range("M17").Value = 57.746478873239
$ws.Range("N17").Value = 57.746478873239

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 81
$ws.Range("J18").Value = 51
$ws.Range("K18").Value = 58.823529411764
$ws.Range("L18").Value = 19.117647058823
$ws.Range("M18").Value = -10
$ws.Range("N18").Value = -82.391304347826

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -38.461538461538
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 68
$ws.Range("H19").Value = -23.529411764705
$ws.Range("I19").Value = 204
$ws.Range("J19").Value = 278
$ws.Range("K19").Value = -26.618705035971
$ws.Range("L19").Value = 22.155688622754
$ws.Range("M19").Value = 129.213483146067
$ws.Range("N19").Value = 36.912751677852

# Row 20
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -40
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 34
$ws.Range("H20").Value = -11.764705882352
$ws.Range("I20").Value = 138
$ws.Range("J20").Value = 134
$ws.Range("K20").Value = 2.985074626865
$ws.Range("L20").Value = 14.049586776859
$ws.Range("M20").Value = 115.625
$ws.Range("N20").Value = -71.721311475409

# Row 21
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -26.470588235294
$ws.Range("F21").Value = 143
$ws.Range("G21").Value = 170
$ws.Range("H21").Value = -15.882352941176
$ws.Range("I21").Value = 614
$ws.Range("J21").Value = 671
$ws.Range("K21").Value = -8.494783904619
$ws.Range("L21").Value = 14.338919925512
$ws.Range("M21").Value = 53.117206982543
$ws.Range("N21").Value = -55.922469490308

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("M22").Value = -42.857142857142

# Row 23
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 9
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 28.571428571428
$ws.Range("I23").Value = 31
$ws.Range("J23").Value = 34
$ws.Range("K23").Value = -8.823529411764
$ws.Range("L23").Value = -18.421052631578
$ws.Range("M23").Value = 72.222222222222

# Row 24
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 81.818181818181
$ws.Range("F24").Value = 135
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = 25
$ws.Range("I24").Value = 406
$ws.Range("J24").Value = 442
$ws.Range("K24").Value = -8.144796380090
$ws.Range("L24").Value = -5.140186915887
$ws.Range("M24").Value = 69.166666666666

# Row 25
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = -8.695652173913
$ws.Range("I25").Value = 132
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = -34
$ws.Range("L25").Value = -26.666666666666

# Row 26
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 45
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = 15.384615384615
$ws.Range("I26").Value = 181
$ws.Range("J26").Value = 143
$ws.Range("K26").Value = 26.573426573426
$ws.Range("L26").Value = 17.532467532467
$ws.Range("M26").Value = 3.428571428571

# Row 27
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -85.714285714285
$ws.Range("J27").Value = 14
$ws.Range("K27").Value = -35.714285714285
$ws.Range("L27").Value = -10

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -80
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -25
$ws.Range("I28").Value = 21
$ws.Range("J28").Value = 21
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 31.25

# Row 29
$ws.Range("N29").Value = -76.923076923076

# Row 30
$ws.Range("N30").Value = -80

# Row 31
$ws.Range("L31").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L31").Value = 0

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "0"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "***.*"

Write-Host "edit applied"
